# Auto-generated edit script updating profit-calculation columns (H:N)
# across multiple sheets per the commit diff (scheduled price-data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 826.7222
$ws.Cells.Item(19, 9).Value = 767.8
$ws.Cells.Item(19, 10).Value = 849.38464
$ws.Cells.Item(19, 11).Value = 767.8
$ws.Cells.Item(19, 12).Value = 849.38464
$ws.Cells.Item(19, 13).Value = -592.8
$ws.Cells.Item(19, 14).Value = -1199.38464
$ws.Cells.Item(41, 8).Value = 623.2632
$ws.Cells.Item(41, 9).Value = 325
$ws.Cells.Item(41, 10).Value = 840.1818
$ws.Cells.Item(41, 11).Value = 325
$ws.Cells.Item(41, 12).Value = 840.1818
$ws.Cells.Item(41, 13).Value = 115
$ws.Cells.Item(41, 14).Value = -1720.1818
$ws.Cells.Item(62, 8).Value = 29025.46
$ws.Cells.Item(62, 9).Value = 6594.4443
$ws.Cells.Item(62, 10).Value = 79495.25
$ws.Cells.Item(62, 11).Value = 6594.4443
$ws.Cells.Item(62, 12).Value = 79495.25
$ws.Cells.Item(62, 13).Value = -5970.4443
$ws.Cells.Item(62, 14).Value = -80743.25
$ws.Cells.Item(65, 8).Value = 29025.46
$ws.Cells.Item(65, 9).Value = 6594.4443
$ws.Cells.Item(65, 10).Value = 79495.25
$ws.Cells.Item(65, 11).Value = 32972.2215
$ws.Cells.Item(65, 12).Value = 397476.25
$ws.Cells.Item(65, 13).Value = -29852.2215
$ws.Cells.Item(65, 14).Value = -403716.25
$ws.Cells.Item(129, 8).Value = 1534.5
$ws.Cells.Item(129, 9).Value = 958.8
$ws.Cells.Item(129, 10).Value = 1726.4
$ws.Cells.Item(129, 11).Value = 2876.4
$ws.Cells.Item(129, 12).Value = 5179.200000000001
$ws.Cells.Item(129, 13).Value = 2123.6
$ws.Cells.Item(129, 14).Value = -15179.2
$ws.Cells.Item(137, 8).Value = 1950
$ws.Cells.Item(137, 9).Value = 2025
$ws.Cells.Item(137, 10).Value = 1821.4286
$ws.Cells.Item(137, 11).Value = 6075
$ws.Cells.Item(137, 12).Value = 5464.2858
$ws.Cells.Item(137, 13).Value = -3525
$ws.Cells.Item(137, 14).Value = -10564.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2526.5334
$ws.Cells.Item(2, 9).Value = 2299.8462
$ws.Cells.Item(2, 10).Value = 4000
$ws.Cells.Item(2, 11).Value = 2299.8462
$ws.Cells.Item(2, 12).Value = 4000
$ws.Cells.Item(2, 13).Value = -2186.8462
$ws.Cells.Item(2, 14).Value = -4226
$ws.Cells.Item(74, 8).Value = 1710.9269
$ws.Cells.Item(74, 9).Value = 1562.6666
$ws.Cells.Item(74, 10).Value = 1996.8572
$ws.Cells.Item(74, 11).Value = 1562.6666
$ws.Cells.Item(74, 12).Value = 1996.8572
$ws.Cells.Item(74, 13).Value = -688.6666
$ws.Cells.Item(74, 14).Value = -3744.8572
$ws.Cells.Item(77, 8).Value = 1710.9269
$ws.Cells.Item(77, 9).Value = 1562.6666
$ws.Cells.Item(77, 10).Value = 1996.8572
$ws.Cells.Item(77, 11).Value = 7813.333000000001
$ws.Cells.Item(77, 12).Value = 9984.286
$ws.Cells.Item(77, 13).Value = -3445.333000000001
$ws.Cells.Item(77, 14).Value = -18720.286
$ws.Cells.Item(97, 8).Value = 1319
$ws.Cells.Item(97, 9).Value = 1157.8572
$ws.Cells.Item(97, 10).Value = 1770.2
$ws.Cells.Item(97, 11).Value = 1157.8572
$ws.Cells.Item(97, 12).Value = 1770.2
$ws.Cells.Item(97, 13).Value = -661.8571999999999
$ws.Cells.Item(97, 14).Value = -2762.2
$ws.Cells.Item(107, 8).Value = 29000
$ws.Cells.Item(107, 10).Value = 29000
$ws.Cells.Item(107, 12).Value = 29000
$ws.Cells.Item(107, 14).Value = -36680
$ws.Cells.Item(109, 8).Value = 14263.333
$ws.Cells.Item(109, 10).Value = 14263.333
$ws.Cells.Item(109, 12).Value = 14263.333
$ws.Cells.Item(109, 14).Value = -17037.333
$ws.Cells.Item(116, 8).Value = 2526.5334
$ws.Cells.Item(116, 9).Value = 2299.8462
$ws.Cells.Item(116, 10).Value = 4000
$ws.Cells.Item(116, 11).Value = 2299.8462
$ws.Cells.Item(116, 12).Value = 4000
$ws.Cells.Item(116, 13).Value = -5.846199999999953
$ws.Cells.Item(116, 14).Value = -8588

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2526.5334
$ws.Cells.Item(3, 9).Value = 2299.8462
$ws.Cells.Item(3, 10).Value = 4000
$ws.Cells.Item(3, 11).Value = 2299.8462
$ws.Cells.Item(3, 12).Value = 4000
$ws.Cells.Item(3, 13).Value = -2185.8462
$ws.Cells.Item(3, 14).Value = -4228
$ws.Cells.Item(94, 8).Value = 940.6429000000001
$ws.Cells.Item(94, 9).Value = 767.4167
$ws.Cells.Item(94, 10).Value = 1980
$ws.Cells.Item(94, 11).Value = 767.4167
$ws.Cells.Item(94, 12).Value = 1980
$ws.Cells.Item(94, 13).Value = -316.4167
$ws.Cells.Item(94, 14).Value = -2882
$ws.Cells.Item(134, 8).Value = 4742.7144
$ws.Cells.Item(134, 9).Value = 2114.7778
$ws.Cells.Item(134, 10).Value = 7967.909
$ws.Cells.Item(134, 11).Value = 6344.3334
$ws.Cells.Item(134, 12).Value = 23903.727
$ws.Cells.Item(134, 13).Value = -3809.3334
$ws.Cells.Item(134, 14).Value = -28973.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 3194.9656
$ws.Cells.Item(86, 9).Value = 2765.875
$ws.Cells.Item(86, 10).Value = 3723.077
$ws.Cells.Item(86, 11).Value = 2765.875
$ws.Cells.Item(86, 12).Value = 3723.077
$ws.Cells.Item(86, 13).Value = -1642.875
$ws.Cells.Item(86, 14).Value = -5969.077
$ws.Cells.Item(89, 8).Value = 3194.9656
$ws.Cells.Item(89, 9).Value = 2765.875
$ws.Cells.Item(89, 10).Value = 3723.077
$ws.Cells.Item(89, 11).Value = 13829.375
$ws.Cells.Item(89, 12).Value = 18615.385
$ws.Cells.Item(89, 13).Value = -8213.375
$ws.Cells.Item(89, 14).Value = -29847.385
$ws.Cells.Item(132, 8).Value = 3427.3572
$ws.Cells.Item(132, 9).Value = 2569.6667
$ws.Cells.Item(132, 10).Value = 4971.2
$ws.Cells.Item(132, 11).Value = 7709.000100000001
$ws.Cells.Item(132, 12).Value = 14913.6
$ws.Cells.Item(132, 13).Value = -5179.000100000001
$ws.Cells.Item(132, 14).Value = -19973.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1005.2941
$ws.Cells.Item(97, 9).Value = 1051.9333
$ws.Cells.Item(97, 10).Value = 655.5
$ws.Cells.Item(97, 11).Value = 1051.9333
$ws.Cells.Item(97, 12).Value = 655.5
$ws.Cells.Item(97, 13).Value = -555.9332999999999
$ws.Cells.Item(97, 14).Value = -1647.5
$ws.Cells.Item(107, 8).Value = 634.8823
$ws.Cells.Item(107, 9).Value = 599.53845
$ws.Cells.Item(107, 10).Value = 749.75
$ws.Cells.Item(107, 11).Value = 599.53845
$ws.Cells.Item(107, 12).Value = 749.75
$ws.Cells.Item(107, 13).Value = 1320.46155
$ws.Cells.Item(107, 14).Value = -4589.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2082.8572
$ws.Cells.Item(22, 9).Value = 2890
$ws.Cells.Item(22, 10).Value = 1760
$ws.Cells.Item(22, 11).Value = 2890
$ws.Cells.Item(22, 12).Value = 1760
$ws.Cells.Item(22, 13).Value = -2595
$ws.Cells.Item(22, 14).Value = -2350
$ws.Cells.Item(27, 8).Value = 2082.8572
$ws.Cells.Item(27, 9).Value = 2890
$ws.Cells.Item(27, 10).Value = 1760
$ws.Cells.Item(27, 11).Value = 2890
$ws.Cells.Item(27, 12).Value = 1760
$ws.Cells.Item(27, 13).Value = -2783
$ws.Cells.Item(27, 14).Value = -1974

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 7913.25
$ws.Cells.Item(96, 9).Value = 4074.75
$ws.Cells.Item(96, 10).Value = 9832.5
$ws.Cells.Item(96, 11).Value = 4074.75
$ws.Cells.Item(96, 12).Value = 9832.5
$ws.Cells.Item(96, 13).Value = -2701.75
$ws.Cells.Item(96, 14).Value = -12578.5
$ws.Cells.Item(126, 8).Value = 2756.7715
$ws.Cells.Item(126, 9).Value = 2675.5217
$ws.Cells.Item(126, 10).Value = 2912.5
$ws.Cells.Item(126, 11).Value = 8026.5651
$ws.Cells.Item(126, 12).Value = 8737.5
$ws.Cells.Item(126, 13).Value = -5556.5651
$ws.Cells.Item(126, 14).Value = -13677.5
$ws.Cells.Item(136, 8).Value = 6946187
$ws.Cells.Item(136, 9).Value = 10870485
$ws.Cells.Item(136, 11).Value = 32611455
$ws.Cells.Item(136, 13).Value = -32608905
